$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPPbI")

# New header cells for energy-related vs. process emissions columns
$ws.Range("B1").Value = "energy related emissions"
$ws.Range("C1").Value = "process emissions"

# Mirror column B values into new column C for each industry row
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 1

# Column width updates (A widened, B/C newly sized for the emissions-split columns).
# Values are chosen so the engine's pixel-quantized ColumnWidth setter lands as
# close as possible to the target stored widths (43.7109375 / 24.7109375 / 25.85546875).
$ws.Columns.Item(1).ColumnWidth = 42.8333333333
$ws.Columns.Item(2).ColumnWidth = 23.8333333333
$ws.Columns.Item(3).ColumnWidth = 25

